$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("BM")

# Fill in previously-blank rows 9 (gamma=0.4) and 10 (gamma=0.3) with measured values
$ws.Range("B9").Value = 0.19850000000000001
$ws.Range("C9").Value = 0.1986
$ws.Range("D9").Value = 0.14460000000000001
$ws.Range("E9").Value = 0.14410000000000001
$ws.Range("F9").Value = 0.1991
$ws.Range("G9").Value = 0.14430000000000001
$ws.Range("H9").Value = 0.1971

$ws.Range("B10").Value = 0.19769999999999999
$ws.Range("C10").Value = 0.19919999999999999
$ws.Range("D10").Value = 0.15279999999999999
$ws.Range("E10").Value = 0.15329999999999999
$ws.Range("F10").Value = 0.19639999999999999
$ws.Range("G10").Value = 0.15279999999999999
$ws.Range("H10").Value = 0.1973

# Update the view state to match: scrolled so row10 is the top visible row, selection on I25
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("I25").Select()
